$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Simple Test"
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
